$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 229: A and C were stored as text, convert to numeric
$ws.Cells.Item(229, 1).Value = 4.26
$ws.Cells.Item(229, 3).Value = 4833

# Append new rows 230-279 (A, C numeric; B text timestamp)
$data = @(
  ,(8.279999999999999, "2024-08-06 13:31:57", 4837)
  ,(1.17, "2024-08-06 13:32:19", 5140)
  ,(1.31, "2024-08-06 13:32:40", 4425)
  ,(1.13, "2024-08-06 13:32:58", 4057)
  ,(2.79, "2024-08-06 13:33:27", 3652)
  ,(2.34, "2024-08-06 13:33:56", 4729)
  ,(1.86, "2024-08-06 13:34:28", 5047)
  ,(4.69, "2024-08-06 13:35:02", 4347)
  ,(1.11, "2024-08-06 13:35:22", 4910)
  ,(1.1, "2024-08-06 13:35:38", 4572)
  ,(1.01, "2024-08-06 13:35:53", 3593)
  ,(1.21, "2024-08-06 13:36:12", 3226)
  ,(4.96, "2024-08-06 13:36:47", 2874)
  ,(1.15, "2024-08-06 13:37:11", 3969)
  ,(1.02, "2024-08-06 13:37:28", 3568)
  ,(1.08, "2024-08-06 13:37:46", 2801)
  ,(3.7, "2024-08-06 13:38:14", 2590)
  ,(1.34, "2024-08-06 13:38:32", 3370)
  ,(3.08, "2024-08-06 13:39:00", 2774)
  ,(2.53, "2024-08-06 13:39:25", 2992)
  ,(2.83, "2024-08-06 13:39:51", 3383)
  ,(1.42, "2024-08-06 13:40:11", 3570)
  ,(1.24, "2024-08-06 13:40:31", 2987)
  ,(1.07, "2024-08-06 13:40:48", 2654)
  ,(1.26, "2024-08-06 13:41:07", 2410)
  ,(1.22, "2024-08-06 13:41:26", 2347)
  ,(2.89, "2024-08-06 13:41:55", 2251)
  ,(1.48, "2024-08-06 13:42:18", 2892)
  ,(1.16, "2024-08-06 13:42:35", 2635)
  ,(1.45, "2024-08-06 13:42:52", 2023)
  ,(1.3, "2024-08-06 13:43:13", 1926)
  ,(4.99, "2024-08-06 13:43:44", 1921)
  ,(6.51, "2024-08-06 13:44:20", 2616)
  ,(1.66, "2024-08-06 13:44:41", 3040)
  ,(1.15, "2024-08-06 13:45:01", 2674)
  ,(2.74, "2024-08-06 13:45:29", 2262)
  ,(1.05, "2024-08-06 13:45:44", 2980)
  ,(1.84, "2024-08-06 13:46:10", 2295)
  ,(2.16, "2024-08-06 13:46:31", 2503)
  ,(1.99, "2024-08-06 13:46:54", 2374)
  ,(44.77, "2024-08-06 13:47:51", 2267)
  ,(3.37, "2024-08-06 13:48:18", 2820)
  ,(3.07, "2024-08-06 13:48:46", 3059)
  ,(26.65, "2024-08-06 13:49:38", 3321)
  ,(2.94, "2024-08-06 13:50:07", 3686)
  ,(1.31, "2024-08-06 13:50:26", 3744)
  ,(1.26, "2024-08-06 13:50:49", 3290)
  ,(1.25, "2024-08-06 13:51:06", 3006)
  ,(6.59, "2024-08-06 13:51:42", 2941)
  ,(2.3, "2024-08-06 13:52:10", 3988)
)

$r = 230
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}

# Row 280: stored as text (A and C look numeric but are text)
$ws.Cells.Item(280, 1).Value = "'1.08"
$ws.Cells.Item(280, 2).Value = "2024-08-06 13:52:31"
$ws.Cells.Item(280, 3).Value = "'4064"
